$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column "12-nov" right before the
# existing "01-oct." column (column DO), shifting everything from DO onward
# one column to the right (DO -> DP, DP -> DQ, ... ES -> ET).
$spot = $wb.Worksheets.Item("Prix Spot")

$spot.Range("DO1").EntireColumn.Insert()

# New header cell for the inserted date column.
$spot.Range("DO1").Value = "12-nov"

# New data column is empty for every hour row; the source marks missing
# values with a literal dash.
$spot.Range("DO2:DO25").Value = "-"

# --- Sheet "Gaz": append one more day of data (row 148).
$gaz = $wb.Worksheets.Item("Gaz")
$gaz.Range("A148").Formula = "=""2025-11-10"""
$gaz.Range("A148").Copy()
$gaz.Range("A148").PasteSpecial(-4163)
$gaz.Range("B148").Value = 28.925

# --- Sheet "CO2": append the same new day of data (row 148).
$co2 = $wb.Worksheets.Item("CO2")
$co2.Range("A148").Formula = "=""2025-11-10"""
$co2.Range("A148").Copy()
$co2.Range("A148").PasteSpecial(-4163)
$co2.Range("B148").Value = 79.88
